$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

$ws.Range("AC2").Value = 0.98392282958199362
$ws.Range("AD2").Value = 1

$ws.Range("AC3").Value = 0.98447204968944102
$ws.Range("AD3").Value = 0.92150170648464158

$ws.Range("A1:AD3").Select() | Out-Null
